# Applies the "Added automatic coloring to plots" revision to chi_squared.xlsx
#
# Actual observable changes in this revision:
#  1. Sheet2 ("Transposed"): header cells B1/C1 swap order (Dead, Alive instead of
#     Alive, Dead) and a new "For SPSS" long-format (Endotype / Status / Frequency)
#     table is appended in rows 14-32, built from the existing cross-tab numbers.
#  2. Sheet2 gains an explicit portrait page setup.
#  3. Sheet1's first chart ("Chart 6") is repositioned (slid down/left) while
#     keeping the same size.
#  4. Selections/active cells on both sheets changed, and Sheet1's frozen
#     top-left-cell scroll position was reset.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: swap the Dead/Alive header cells -----------------------------
$ws2.Range("B1").Value = "Dead"
$ws2.Range("C1").Value = "Alive"

# --- Sheet2: append the "For SPSS" long-format table ----------------------
$ws2.Range("A14").Value = "For SPSS"
$ws2.Range("C14").Value = "Need to weight by frequency first"

$ws2.Range("A15").Value = "Endotype"
$ws2.Range("B15").Value = "Status"
$ws2.Range("C15").Value = "Frequency"

$endotypes = @("Endotype 1", "Endotype 2", "Endotype 3", "Endotype 4", "Endotype 5", "Endotype 6", "Endotype 7", "Endotype 8")
$deadCounts  = @(47, 441, 147, 33, 121, 320, 76, 54)
$aliveCounts = @(11, 428, 63, 18, 115, 948, 47, 112)

# Dead rows: 16-23
for ($i = 0; $i -lt 8; $i++) {
    $r = 16 + $i
    $ws2.Range("A$r").Value = $endotypes[$i]
    $ws2.Range("B$r").Value = "Dead"
    $ws2.Range("C$r").Value = $deadCounts[$i]
}

# Alive rows: 24-31
for ($i = 0; $i -lt 8; $i++) {
    $r = 24 + $i
    $ws2.Range("A$r").Value = $endotypes[$i]
    $ws2.Range("B$r").Value = "Alive"
    $ws2.Range("C$r").Value = $aliveCounts[$i]
}

$ws2.Range("A32").Value = 'Notes: Weight by frequency (under data→weight), crosstabs, select options for Chi square under "tests", select column comparisons with bonferroni adjustment under "cells", put endotype as columns and status as row'

# --- Sheet2: page setup ----------------------------------------------------
$ws2.PageSetup.Orientation = 1

# --- Sheet1: reposition the first chart (same size, new top-left) --------
$chart = $ws1.ChartObjects().Item(1)
$chart.Left = 139.575
$chart.Top = 233.3

# --- Selections / view state ----------------------------------------------
$ws2.Range("B2").Select()
$ws1.Range("B13").Select()

$wb.Save()
